$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.502.09"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.251.26"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'308.32"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'95.17"
$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").Value = "'35.34"
$ws.Range("E10").Value = "  +3.68%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("D12").Value = "'7.31"
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "2.308.58"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").Value = "'0.841"
$ws.Range("E15").Value = "  +4.15%  "
$ws.Range("D16").Value = "'13.72"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "44.212.08"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "0.0₃0968"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  +5.48%  "
$ws.Range("D20").Value = "'12.24"
$ws.Range("E20").Value = "  +2.65%  "
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("D22").Value = "'238.23"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("E23").Value = "  +3.86%  "
$ws.Range("E24").Value = "  +5.34%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  +5.90%  "
$ws.Range("D27").Value = "'9.90"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").Value = "'38.06"
$ws.Range("E28").Value = "  +5.87%  "
$ws.Range("D29").Value = "'6.01"
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").Value = "'20.12"
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("D31").Value = "'152.68"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'0.0804"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'2.64"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").Value = "'3.19"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "'0.121"
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("E38").Value = "  +6.10%  "
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("D40").Value = "'3.83"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "1.752.07"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +5.93%  "
$ws.Range("D45").Value = "'81.08"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "'100.40"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("D47").Value = "'71.19"
$ws.Range("E47").Value = "  +5.13%  "
$ws.Range("D48").Value = "'55.64"
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("D49").Value = "'8.19"
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("D50").Value = "'1.59"
$ws.Range("E50").Value = "  +7.34%  "
$ws.Range("D51").Value = "'4.86"
$ws.Range("E51").Value = "  -0.47%  "
